$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("AK2").Value = 34
$ws.Range("AS2").Value = 301
$ws.Range("AW2").Value = 5
$ws.Range("AX2").Value = 21
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.63
$ws.Range("K3").Value = 1.91
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 4.75
$ws.Range("AC3").Value = 6
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 21
$ws.Range("AI3").Value = 17
$ws.Range("AL3").Value = 51
$ws.Range("AP3").Value = 29
$ws.Range("AS3").Value = 351
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 10
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 126
$ws.Range("BA3").Value = 201
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("BD5").Value = 126
$ws.Range("G6").Value = 2.05
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 2.6
$ws.Range("Q6").Value = 1.62
$ws.Range("R6").Value = 2.25
$ws.Range("Z6").Value = 19
$ws.Range("AP6").Value = 17
$ws.Range("AR6").Value = 41
$ws.Range("BA6").Value = 67
$ws.Range("BC6").Value = 351
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("G9").Value = 1.53
$ws.Range("K9").Value = 2.1
$ws.Range("AD9").Value = 7.5
$ws.Range("AO9").Value = 8
$ws.Range("O14").Value = 1.4
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.25
$ws.Range("R14").Value = 1.62
$ws.Range("Q15").Value = 1.98
$ws.Range("R15").Value = 1.92
